$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert GAP rows, working from the bottom of the sheet upward so that
# earlier (lower row-number) insertions don't invalidate later row indices.
# Each entry: row index (in the *original* 1-27 numbering, applied bottom-up)
# at which a new blank row is inserted (pushing that row and below down by one),
# followed by the A/B/C/D values to put into the newly created row and the
# style index to apply (3 = yellow "GAP" fill style).

$insertions = @(
    @{ Row = 27; A = "Paleocene"; B = "GAP"; C = 64;   D = 65 },
    @{ Row = 26; A = "Paleocene"; B = "GAP"; C = 61;   D = 63 },
    @{ Row = 25; A = "Paleocene-Eocene"; B = "GAP"; C = 55; D = 59.2 },
    @{ Row = 23; A = "Eocene"; B = "GAP"; C = 48;   D = 50 },
    @{ Row = 22; A = "Eocene"; B = "GAP"; C = 45;   D = 47 },
    @{ Row = 21; A = "Eocene"; B = "GAP"; C = 42;   D = 44 },
    @{ Row = 19; A = "Eocene"; B = "GAP"; C = 33.9; D = 36 },
    @{ Row = 17; A = "Miocene"; B = "GAP"; C = 21;  D = 23.03 },
    @{ Row = 16; A = "Miocene"; B = "GAP"; C = 18;  D = 20 }
)

foreach ($ins in $insertions) {
    $r = $ins.Row
    $ws.Rows.Item($r).Insert()
    $ws.Range("A$r").Value2 = $ins.A
    $ws.Range("B$r").Value2 = $ins.B
    $ws.Range("C$r").Value2 = $ins.C
    $ws.Range("D$r").Value2 = $ins.D
    $ws.Range("A$r`:D$r").Interior.Color = 65535
}

# Update the sheet view to match the post-edit selection state.
$ws.Range("A1:D36").Select()
